$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column AE (column 31)
$ws.Cells.Item(1, 31).Value = "17-jul"

# Data values for AE2:AE18
$values = @(
    0,
    11.916344324914331,
    15.595549801872435,
    26.508480011935593,
    0,
    1.2043022267373138,
    11.549025325762834,
    24.067907073703363,
    24.161856483336027,
    11.990364530593782,
    0,
    8.9142744139749173,
    0,
    0,
    14.483721004046092,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 31).Value = $values[$i]
}

# Update selection to mirror final saved state
$ws.Range("AG9").Select()
